# Add a new row of custom employee data to the "PIM" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PIM")

$ws.Range("A5").Value = "Shreyas"
$ws.Range("B5").Value = "K"
$ws.Range("C5").Value = "Iyer"

$ws.Activate()
$ws.Range("D8").Select()
